$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.564.69"
$ws.Range("E2").Value = "  -1.41%  "
$ws.Range("D3").Value = "2.902.24"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "526.86"
$ws.Range("E5").Value = "  -2.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.87"
$ws.Range("E6").Value = "  -4.18%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -3.12%  "
$ws.Range("D9").Value = "2.908.64"
$ws.Range("E9").Value = "  -2.28%  "
$ws.Range("E10").Value = "  -5.04%  "
$ws.Range("E11").Value = "  -2.57%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.358"
$ws.Range("E12").Value = "  -2.52%  "
$ws.Range("D13").Value = "3.412.07"
$ws.Range("E13").Value = "  -2.09%  "
$ws.Range("E14").Value = "  +2.38%  "
$ws.Range("D15").Value = "60.544.11"
$ws.Range("E15").Value = "  -1.59%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.59"
$ws.Range("E16").Value = "  -4.08%  "
$ws.Range("D17").Value = "2.909.49"
$ws.Range("E17").Value = "  -2.03%  "
$ws.Range("E18").Value = "  -3.57%  "
$ws.Range("E19").Value = "  -3.45%  "
$ws.Range("E20").Value = "  -3.35%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "352.17"
$ws.Range("E21").Value = "  -6.98%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.56"
$ws.Range("E22").Value = "  -1.73%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.70"
$ws.Range("E24").Value = "  +0.79%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "64.62"
$ws.Range("E25").Value = "  -1.40%  "
$ws.Range("E26").Value = "  -3.58%  "
$ws.Range("E27").Value = "  -4.92%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.51%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.81"
$ws.Range("E29").Value = "  -4.77%  "
$ws.Range("D30").Value = "0.0₃0848"
$ws.Range("E30").Value = "  -9.21%  "
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("E32").Value = "  -2.65%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.60"
$ws.Range("E33").Value = "  -3.76%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "151.11"
$ws.Range("E34").Value = "  -5.27%  "
$ws.Range("E35").Value = "  -6.79%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.57"
$ws.Range("E36").Value = "  -5.67%  "
$ws.Range("E37").Value = "  -6.66%  "
$ws.Range("E38").Value = "  -5.23%  "
$ws.Range("E39").Value = "  +0.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.47"
$ws.Range("E40").Value = "  -4.92%  "
$ws.Range("E41").Value = "  -4.98%  "
$ws.Range("D42").Value = "2.290.63"
$ws.Range("E42").Value = "  -5.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.648"
$ws.Range("E43").Value = "  -2.98%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0581"
$ws.Range("E44").Value = "  -1.46%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "20.37"
$ws.Range("E45").Value = "  -7.57%  "
$ws.Range("E46").Value = "  +0.07%  "
$ws.Range("E47").Value = "  -1.53%  "
$ws.Range("E48").Value = "  -3.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0920"
$ws.Range("E50").Value = "  -3.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "247.57"
$ws.Range("E51").Value = "  -7.03%  "
